# fix: time unit corrected
# The report used "clock cycles" instead of proper time units (us/ms)
# in two spots describing the LBIST total overhead / test-phase duration.

$d = $word.ActiveDocument

# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,
#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,
#              Format, ReplaceWith, Replace)
# wdReplaceAll = 2, wdFindContinue = 1

# 1) "... 5050 * 12 * 10 = 606 K clock cycles." -> "... = 606 us."
$d.Content.Find.Execute(
    "606 K clock cycles.", $true, $false, $false, $false, $false,
    $true, 1, $false, "606 us.", 2) | Out-Null

# 2) "...needs at least 0.6 M clock cycles during which 60 K patterns..."
#    -> "...needs at least 0.6 ms during which 60 K patterns..."
$d.Content.Find.Execute(
    "0.6 M clock cycles during which 60", $true, $false, $false, $false, $false,
    $true, 1, $false, "0.6 ms during which 60", 2) | Out-Null
